$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 633265.1
$ws.Range("J17").Value = 662474.3
$ws.Range("L17").Value = 1987422.9
$ws.Range("N17").Value = -1987758.9

$ws.Range("H32").Value = 1421.3334
$ws.Range("J32").Value = 1632
$ws.Range("L32").Value = 1632
$ws.Range("N32").Value = -2284

$ws.Range("H62").Value = 9241.026
$ws.Range("I62").Value = 6632.12
$ws.Range("J62").Value = 14258.154
$ws.Range("K62").Value = 6632.12
$ws.Range("L62").Value = 14258.154
$ws.Range("M62").Value = -6008.12
$ws.Range("N62").Value = -15506.154

$ws.Range("H65").Value = 9241.026
$ws.Range("I65").Value = 6632.12
$ws.Range("J65").Value = 14258.154
$ws.Range("K65").Value = 33160.6
$ws.Range("L65").Value = 71290.77
$ws.Range("M65").Value = -30040.6
$ws.Range("N65").Value = -77530.77

$ws.Range("H76").Value = 3587300.2
$ws.Range("J76").Value = 3783.3333
$ws.Range("L76").Value = 3783.3333
$ws.Range("N76").Value = -4413.3333

$ws.Range("H79").Value = 3587300.2
$ws.Range("J79").Value = 3783.3333
$ws.Range("L79").Value = 3783.3333
$ws.Range("N79").Value = -5967.3333

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H86").Value = 1375.75
$ws.Range("I86").Value = 802
$ws.Range("J86").Value = 1720
$ws.Range("K86").Value = 802
$ws.Range("L86").Value = 1720
$ws.Range("M86").Value = 321
$ws.Range("N86").Value = -3966

$ws.Range("H89").Value = 1375.75
$ws.Range("I89").Value = 802
$ws.Range("J89").Value = 1720
$ws.Range("K89").Value = 4010
$ws.Range("L89").Value = 8600
$ws.Range("M89").Value = 1606
$ws.Range("N89").Value = -19832

$ws.Range("H129").Value = 2128.625
$ws.Range("I129").Value = 535.25
$ws.Range("J129").Value = 2659.75
$ws.Range("K129").Value = 1605.75
$ws.Range("L129").Value = 7979.25
$ws.Range("M129").Value = 3394.25
$ws.Range("N129").Value = -17979.25

$ws.Range("H133").Value = 23103.75
$ws.Range("J133").Value = 23103.75
$ws.Range("L133").Value = 23103.75
$ws.Range("N133").Value = -33223.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21349.686
$ws.Range("I32").Value = 2824.4038
$ws.Range("K32").Value = 2824.4038
$ws.Range("M32").Value = -2537.4038

$ws.Range("H122").Value = 2446.5454
$ws.Range("I122").Value = 2358.4211
$ws.Range("K122").Value = 7075.263300000001
$ws.Range("M122").Value = -4625.263300000001

$ws.Range("H139").Value = 30729.5
$ws.Range("J139").Value = 44859
$ws.Range("L139").Value = 44859
$ws.Range("N139").Value = -55139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1143.8
$ws.Range("I94").Value = 889.8
$ws.Range("K94").Value = 889.8
$ws.Range("M94").Value = -438.8

$ws.Range("H99").Value = 1830.04
$ws.Range("I99").Value = 1592
$ws.Range("K99").Value = 1592
$ws.Range("M99").Value = -94

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H108").Value = 29250
$ws.Range("J108").Value = 29250
$ws.Range("L108").Value = 29250
$ws.Range("N108").Value = -36930

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2148.7856
$ws.Range("I6").Value = 3170.25
$ws.Range("J6").Value = 1740.2
$ws.Range("K6").Value = 3170.25
$ws.Range("L6").Value = 1740.2
$ws.Range("M6").Value = -3057.25
$ws.Range("N6").Value = -1966.2

$ws.Range("H31").Value = 1865.1052
$ws.Range("I31").Value = 989.5454999999999
$ws.Range("J31").Value = 3069
$ws.Range("K31").Value = 989.5454999999999
$ws.Range("L31").Value = 3069
$ws.Range("M31").Value = -694.5454999999999
$ws.Range("N31").Value = -3659

$ws.Range("H34").Value = 1865.1052
$ws.Range("I34").Value = 989.5454999999999
$ws.Range("J34").Value = 3069
$ws.Range("K34").Value = 989.5454999999999
$ws.Range("L34").Value = 3069
$ws.Range("M34").Value = -787.5454999999999
$ws.Range("N34").Value = -3473

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 483.1905
$ws.Range("I11").Value = 560.5833
$ws.Range("J11").Value = 380
$ws.Range("K11").Value = 1681.7499
$ws.Range("L11").Value = 1140
$ws.Range("M11").Value = -1541.7499
$ws.Range("N11").Value = -1420

$ws.Range("H94").Value = 2526.2632
$ws.Range("J94").Value = 2973.3333
$ws.Range("L94").Value = 8919.999899999999
$ws.Range("N94").Value = -10271.9999

$ws.Range("H131").Value = 1830.5758
$ws.Range("J131").Value = 2217.2693
$ws.Range("L131").Value = 6651.8079
$ws.Range("N131").Value = -16731.8079

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2222.6538
$ws.Range("I126").Value = 1769.8572
$ws.Range("J126").Value = 2389.4736
$ws.Range("K126").Value = 5309.571599999999
$ws.Range("L126").Value = 7168.4208
$ws.Range("M126").Value = -2839.571599999999
$ws.Range("N126").Value = -12108.4208

$ws.Range("H138").Value = 63350
$ws.Range("J138").Value = 63350
$ws.Range("L138").Value = 63350
$ws.Range("N138").Value = -73630

$ws.Range("H139").Value = 34000
$ws.Range("J139").Value = 34000
$ws.Range("L139").Value = 34000
$ws.Range("N139").Value = -44280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 6336.6665
$ws.Range("I9").Value = 6005
$ws.Range("J9").Value = 7000
$ws.Range("K9").Value = 6005
$ws.Range("L9").Value = 7000
$ws.Range("M9").Value = -5781
$ws.Range("N9").Value = -7448

$ws.Range("H55").Value = 264.73334
$ws.Range("I55").Value = 232.08333
$ws.Range("J55").Value = 286.5
$ws.Range("K55").Value = 232.08333
$ws.Range("L55").Value = 286.5
$ws.Range("M55").Value = -59.08332999999999
$ws.Range("N55").Value = -632.5

$ws.Range("H68").Value = 2481.4707
$ws.Range("I68").Value = 2243.7144
$ws.Range("J68").Value = 2647.9
$ws.Range("K68").Value = 2243.7144
$ws.Range("L68").Value = 2647.9
$ws.Range("M68").Value = -1494.7144
$ws.Range("N68").Value = -4145.9

$ws.Range("H71").Value = 2481.4707
$ws.Range("I71").Value = 2243.7144
$ws.Range("J71").Value = 2647.9
$ws.Range("K71").Value = 11218.572
$ws.Range("L71").Value = 13239.5
$ws.Range("M71").Value = -7474.572
$ws.Range("N71").Value = -20727.5

$ws.Range("H82").Value = 911.087
$ws.Range("I82").Value = 775.2
$ws.Range("J82").Value = 1015.61536
$ws.Range("K82").Value = 775.2
$ws.Range("L82").Value = 1015.61536
$ws.Range("M82").Value = -414.2
$ws.Range("N82").Value = -1737.61536

$ws.Range("H85").Value = 911.087
$ws.Range("I85").Value = 775.2
$ws.Range("J85").Value = 1015.61536
$ws.Range("K85").Value = 775.2
$ws.Range("L85").Value = 1015.61536
$ws.Range("M85").Value = 472.8
$ws.Range("N85").Value = -3511.61536

$ws.Range("H94").Value = 19990
$ws.Range("J94").Value = 19990
$ws.Range("L94").Value = 19990
$ws.Range("N94").Value = -21342

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 21931.389
$ws.Range("I14").Value = 16800
$ws.Range("J14").Value = 22957.666
$ws.Range("K14").Value = 16800
$ws.Range("L14").Value = 22957.666
$ws.Range("M14").Value = -16632
$ws.Range("N14").Value = -23293.666

$ws.Range("H132").Value = 10418997
$ws.Range("I132").Value = 15627097
$ws.Range("K132").Value = 46881291
$ws.Range("M132").Value = -46878761
